# First report to COLBE done
# Adds a new "Comment" column (F) with a handful of comments on specific
# parameter rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column F.
$ws.Range("F1").Value = "Comment"

# Row-specific comments (order matches the shared-string insertion order
# observed in the target workbook).
$ws.Range("F16").Value = "Number of internal paritions"
$ws.Range("F18").Value = "Heat in Watts"
$ws.Range("F24").Value = "Litres per person per sec?"
$ws.Range("F31").Value = "Types (options 1, 2, 3?)"
$ws.Range("F32").Value = "Thermal mass of external walls"
$ws.Range("F33").Value = "Roof"
$ws.Range("F2").Value = "Alfonso (type of building)"

# Leave the cursor where the author left it.
$ws.Range("F9").Select()
